# Add two new UI-locator entries (search-by-employee-name feature) to the
# "Web" worksheet's locator table, following the existing ElementID /
# ElementPath / Method layout used throughout the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22: locator for the "search by employee name" input field
$ws.Cells.Item(22, 1).Value = "UM_searchEmployeeName"
$ws.Cells.Item(22, 2).Value = "//*[@id=""app""]/div[1]/div[2]/div[2]/div/div[1]/div[2]/form/div[1]/div/div[3]/div/div[2]/div/div/input"
$ws.Cells.Item(22, 3).Value = "By.xpath"

# Row 23: locator for the employee-name search result
$ws.Cells.Item(23, 1).Value = "UM_searchResultEmployeeName"
$ws.Cells.Item(23, 2).Value = "//*[@id=""app""]/div[1]/div[2]/div[2]/div/div[2]/div[3]/div/div[2]/div/div/div[4]/div"
$ws.Cells.Item(23, 3).Value = "By.xpath"

# Rows 24-31: blank spacer rows (mirrors the existing blank row 5 pattern
# used elsewhere in the sheet to separate locator groups).
for ($r = 24; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = ""
    $ws.Cells.Item($r, 2).Value = ""
    $ws.Cells.Item($r, 3).Value = ""
}
